$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.460.12"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "'1.872.88"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D5").Value = "'313.56"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").Value = "'0.4794"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").Value = "'0.3761"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("D9").Value = "'0.07380"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").Value = "'0.9416"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").Value = "'20.71"
$ws.Range("D12").Value = "'0.07894"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("D13").Value = "'1.893.59"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "'5.437"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "'6.607"
$ws.Range("D16").Value = "'90.88"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "'1.017"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'0.000008924"
$ws.Range("E18").Value = "  +3.21%  "
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'14.92"
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").Value = "'27.482.36"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").Value = "'5.147"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'1.955"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'154.26"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").Value = "'18.58"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").Value = "'2.018"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'116.15"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").Value = "'5.013"
$ws.Range("E29").Value = "  +3.07%  "
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "'1.218"
$ws.Range("E32").Value = "  +4.40%  "
$ws.Range("D33").Value = "'4.592"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").Value = "'0.7483"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").Value = "'2.695"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("D38").Value = "'0.05301"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").Value = "'3.001"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").Value = "'0.5368"
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("D41").Value = "'7.081"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").Value = "'0.1529"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "'8.424"
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.4845"
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.60"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("D47").Value = "'1.665"
$ws.Range("E47").Value = "  +3.93%  "
$ws.Range("D48").Value = "'103.27"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Value = "'67.17"
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("D50").Value = "'0.06105"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("D51").Value = "'0.9013"
$ws.Range("E51").Value = "  +1.96%  "
